$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "  订单编号  " / " 物流公司  " header labels in A1/B1 were padded with stray
# leading/trailing spaces and (for A1) split across two runs/fonts left over from
# manual edits. That caused the "订单编号"/"物流公司" column headers to fail exact
# matches downstream (the "empty unknown_column_name list" warning) since the
# trimmed header name never matched the padded cell text.
#
# Clean them up to plain, trimmed text, and align A1's font with B1's so both
# headers share one consistent bold title style.
$ws.Range("A1").Font.Name = "Microsoft YaHei"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 15
$ws.Range("A1").Font.Color = 0
$ws.Range("A1").Value = "订单编号"
$ws.Range("B1").Value = "物流公司"

# Restore the cursor to the cell it was left on when the fix was made.
[void]$ws.Range("B6").Select()
